$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Melany Vivar Tranquilo -> King Alale ---
$ws.Range("A4").Value = "King Alale"
$ws.Range("A4").Locked = $true

# --- Row 5: Olasunkanmi Akinmuyipitan -> Roksalana (trailing space) ---
$ws.Range("A5").Value = "Roksalana "

# --- Row 6: Raquel Eusébio -> Shinde Gaikwad ---
$ws.Range("A6").Value = "Shinde Gaikwad"

# --- Row 7: Shinde Gaikwad -> Swahna ---
$ws.Range("A7").Value = "Swahna"

# --- Row 8: Monday shift OFF, cell becomes locked ---
$ws.Range("B8").Value = "OFF"
$ws.Range("B8").Locked = $true

# --- Row 9: Monday shift OFF, cell becomes locked ---
$ws.Range("B9").Value = "OFF"
$ws.Range("B9").Locked = $true

# --- Row 11: Monday shift 13.00 pm to 21.00 pm -> 15.00 pm to 23.00 pm ---
$ws.Range("B11").Value = "15.00 pm to 23.00 pm"
$ws.Range("B11").Locked = $true

# --- Row 12: Monday shift 15.00 pm to 23.00 pm -> OFF ---
$ws.Range("B12").Value = "OFF"
$ws.Range("B12").Locked = $true

# --- Row 16: Monday shift OFF -> 13.00 pm to 21.00 pm ---
$ws.Range("B16").Value = "13.00 pm to 21.00 pm"
$ws.Range("B16").Locked = $true

# --- Row 17 (new): Raquel, Monday shift 13.00 pm to 21.00 pm ---
$ws.Range("A17").Value = "Raquel"
$ws.Range("B17").Value = "13.00 pm to 21.00 pm"
$ws.Range("B17").Locked = $true

# --- Scroll the sheet view down by one row (topLeftCell = A2) ---
$ws.Application.ActiveWindow.ScrollRow = 2
